$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(38,2).Value2 = 6782522
$ws.Cells.Item(38,6).Value2 = 'Municipal Perez Zeledon'
$ws.Cells.Item(38,7).Value2 = 'Sporting San Jose'
$ws.Cells.Item(38,9).Value2 = 2
$ws.Cells.Item(38,10).Value2 = 'A'
$ws.Cells.Item(38,11).Value2 = 2.5
$ws.Cells.Item(38,12).Value2 = 3.5
$ws.Cells.Item(38,13).Value2 = 2.5
$ws.Cells.Item(38,14).Value2 = 2.2
$ws.Cells.Item(38,15).Value2 = 3.5
$ws.Cells.Item(38,16).Value2 = 2.9
$ws.Cells.Item(38,18).Value2 = 1.9
$ws.Cells.Item(38,19).Value2 = 1.9
$ws.Cells.Item(38,20).Value2 = 2.5
$ws.Cells.Item(38,23).Value2 = -1
$ws.Cells.Item(38,25).Value2 = 1.9
$ws.Cells.Item(38,26).Value2 = -1
$ws.Cells.Item(38,27).Value2 = 0.8999999999999999
$ws.Cells.Item(38,28).Value2 = 0.8999999999999999
$ws.Cells.Item(38,29).Value2 = -1
$ws.Cells.Item(39,2).Value2 = 6781354
$ws.Cells.Item(39,6).Value2 = 'Puntarenas'
$ws.Cells.Item(39,7).Value2 = 'AD San Carlos'
$ws.Cells.Item(39,9).Value2 = 0
$ws.Cells.Item(39,10).Value2 = 'H'
$ws.Cells.Item(39,11).Value2 = 2.4
$ws.Cells.Item(39,12).Value2 = 3.2
$ws.Cells.Item(39,13).Value2 = 2.8
$ws.Cells.Item(39,14).Value2 = 2.3
$ws.Cells.Item(39,15).Value2 = 3.2
$ws.Cells.Item(39,16).Value2 = 3
$ws.Cells.Item(39,18).Value2 = 2
$ws.Cells.Item(39,19).Value2 = 1.8
$ws.Cells.Item(39,20).Value2 = 2.25
$ws.Cells.Item(39,23).Value2 = 1.3
$ws.Cells.Item(39,25).Value2 = -1
$ws.Cells.Item(39,26).Value2 = 1
$ws.Cells.Item(39,27).Value2 = -1
$ws.Cells.Item(39,28).Value2 = -1
$ws.Cells.Item(39,29).Value2 = 0.8999999999999999
$ws.Cells.Item(91,2).Value2 = 6782566
$ws.Cells.Item(91,6).Value2 = 'Cartagines'
$ws.Cells.Item(91,7).Value2 = 'Deportivo Saprissa'
$ws.Cells.Item(91,8).Value2 = 0
$ws.Cells.Item(91,9).Value2 = 4
$ws.Cells.Item(91,10).Value2 = 'A'
$ws.Cells.Item(91,11).Value2 = 3.2
$ws.Cells.Item(91,12).Value2 = 3.4
$ws.Cells.Item(91,13).Value2 = 2
$ws.Cells.Item(91,14).Value2 = 2.9
$ws.Cells.Item(91,15).Value2 = 3.5
$ws.Cells.Item(91,16).Value2 = 2.15
$ws.Cells.Item(91,17).Value2 = 0.25
$ws.Cells.Item(91,18).Value2 = 1.875
$ws.Cells.Item(91,19).Value2 = 1.925
$ws.Cells.Item(91,20).Value2 = 3
$ws.Cells.Item(91,21).Value2 = 1.975
$ws.Cells.Item(91,22).Value2 = 1.825
$ws.Cells.Item(91,24).Value2 = -1
$ws.Cells.Item(91,25).Value2 = 1.15
$ws.Cells.Item(91,27).Value2 = 0.925
$ws.Cells.Item(91,28).Value2 = 0.9750000000000001
$ws.Cells.Item(91,29).Value2 = -1
$ws.Cells.Item(92,2).Value2 = 6782568
$ws.Cells.Item(92,6).Value2 = 'Sporting San Jose'
$ws.Cells.Item(92,7).Value2 = 'AD Guanacasteca'
$ws.Cells.Item(92,8).Value2 = 1
$ws.Cells.Item(92,9).Value2 = 1
$ws.Cells.Item(92,10).Value2 = 'D'
$ws.Cells.Item(92,11).Value2 = 1.909
$ws.Cells.Item(92,12).Value2 = 3.6
$ws.Cells.Item(92,13).Value2 = 3.3
$ws.Cells.Item(92,14).Value2 = 2
$ws.Cells.Item(92,15).Value2 = 3.6
$ws.Cells.Item(92,16).Value2 = 3.1
$ws.Cells.Item(92,17).Value2 = -0.5
$ws.Cells.Item(92,18).Value2 = 2
$ws.Cells.Item(92,19).Value2 = 1.8
$ws.Cells.Item(92,20).Value2 = 2.5
$ws.Cells.Item(92,21).Value2 = 1.825
$ws.Cells.Item(92,22).Value2 = 1.975
$ws.Cells.Item(92,24).Value2 = 2.6
$ws.Cells.Item(92,25).Value2 = -1
$ws.Cells.Item(92,27).Value2 = 0.8
$ws.Cells.Item(92,28).Value2 = -1
$ws.Cells.Item(92,29).Value2 = 0.9750000000000001
$ws.Cells.Item(110,2).Value2 = 6782579
$ws.Cells.Item(110,6).Value2 = 'Santos de Gupiles'
$ws.Cells.Item(110,7).Value2 = 'AD San Carlos'
$ws.Cells.Item(110,8).Value2 = 0
$ws.Cells.Item(110,9).Value2 = 2
$ws.Cells.Item(110,10).Value2 = 'A'
$ws.Cells.Item(110,11).Value2 = 2.4
$ws.Cells.Item(110,12).Value2 = 3.3
$ws.Cells.Item(110,13).Value2 = 2.7
$ws.Cells.Item(110,14).Value2 = 2.375
$ws.Cells.Item(110,15).Value2 = 3.4
$ws.Cells.Item(110,16).Value2 = 2.8
$ws.Cells.Item(110,17).Value2 = -0.25
$ws.Cells.Item(110,18).Value2 = 2
$ws.Cells.Item(110,19).Value2 = 1.8
$ws.Cells.Item(110,20).Value2 = 2.5
$ws.Cells.Item(110,21).Value2 = 1.875
$ws.Cells.Item(110,22).Value2 = 1.925
$ws.Cells.Item(110,23).Value2 = -1
$ws.Cells.Item(110,25).Value2 = 1.8
$ws.Cells.Item(110,26).Value2 = -1
$ws.Cells.Item(110,27).Value2 = 0.8
$ws.Cells.Item(110,29).Value2 = 0.925
$ws.Cells.Item(111,2).Value2 = 6782581
$ws.Cells.Item(111,6).Value2 = 'Alajuelense'
$ws.Cells.Item(111,7).Value2 = 'AD Grecia'
$ws.Cells.Item(111,8).Value2 = 2
$ws.Cells.Item(111,9).Value2 = 0
$ws.Cells.Item(111,10).Value2 = 'H'
$ws.Cells.Item(111,11).Value2 = 1.181
$ws.Cells.Item(111,12).Value2 = 6.5
$ws.Cells.Item(111,13).Value2 = 11
$ws.Cells.Item(111,14).Value2 = 1.25
$ws.Cells.Item(111,15).Value2 = 5
$ws.Cells.Item(111,16).Value2 = 9
$ws.Cells.Item(111,17).Value2 = -1.75
$ws.Cells.Item(111,18).Value2 = 1.975
$ws.Cells.Item(111,19).Value2 = 1.825
$ws.Cells.Item(111,20).Value2 = 3.25
$ws.Cells.Item(111,21).Value2 = 2
$ws.Cells.Item(111,22).Value2 = 1.8
$ws.Cells.Item(111,23).Value2 = 0.25
$ws.Cells.Item(111,25).Value2 = -1
$ws.Cells.Item(111,26).Value2 = 0.4875
$ws.Cells.Item(111,27).Value2 = -0.5
$ws.Cells.Item(111,29).Value2 = 0.8
$ws.Cells.Item(129,2).Value2 = 6782595
$ws.Cells.Item(129,6).Value2 = 'Herediano'
$ws.Cells.Item(129,7).Value2 = 'Sporting San Jose'
$ws.Cells.Item(129,8).Value2 = 3
$ws.Cells.Item(129,11).Value2 = 1.4
$ws.Cells.Item(129,12).Value2 = 4.75
$ws.Cells.Item(129,13).Value2 = 7
$ws.Cells.Item(129,14).Value2 = 1.363
$ws.Cells.Item(129,15).Value2 = 4.75
$ws.Cells.Item(129,16).Value2 = 8.5
$ws.Cells.Item(129,17).Value2 = -1.25
$ws.Cells.Item(129,20).Value2 = 3
$ws.Cells.Item(129,21).Value2 = 1.95
$ws.Cells.Item(129,22).Value2 = 1.85
$ws.Cells.Item(129,23).Value2 = 0.363
$ws.Cells.Item(129,28).Value2 = 0
$ws.Cells.Item(129,29).Value2 = -0
$ws.Cells.Item(131,2).Value2 = 6782598
$ws.Cells.Item(131,6).Value2 = 'Municipal Perez Zeledon'
$ws.Cells.Item(131,7).Value2 = 'Cartagines'
$ws.Cells.Item(131,8).Value2 = 1
$ws.Cells.Item(131,11).Value2 = 4.5
$ws.Cells.Item(131,12).Value2 = 3.75
$ws.Cells.Item(131,13).Value2 = 1.615
$ws.Cells.Item(131,14).Value2 = 3.4
$ws.Cells.Item(131,15).Value2 = 3.4
$ws.Cells.Item(131,16).Value2 = 1.85
$ws.Cells.Item(131,17).Value2 = 0.5
$ws.Cells.Item(131,20).Value2 = 2.75
$ws.Cells.Item(131,21).Value2 = 1.9
$ws.Cells.Item(131,22).Value2 = 1.9
$ws.Cells.Item(131,23).Value2 = 2.4
$ws.Cells.Item(131,28).Value2 = -1
$ws.Cells.Item(131,29).Value2 = 0.8999999999999999
$ws.Cells.Item(223,2).Value2 = 7980420
$ws.Cells.Item(223,5).Value2 = 45385.95833333334
$ws.Cells.Item(223,6).Value2 = 'Deportivo Saprissa'
$ws.Cells.Item(223,7).Value2 = 'Municipal Perez Zeledon'
$ws.Cells.Item(223,8).Value2 = 1
$ws.Cells.Item(223,9).Value2 = 0
$ws.Cells.Item(223,10).Value2 = 'H'
$ws.Cells.Item(223,11).Value2 = 1.222
$ws.Cells.Item(223,12).Value2 = 6
$ws.Cells.Item(223,13).Value2 = 12
$ws.Cells.Item(223,14).Value2 = 1.166
$ws.Cells.Item(223,15).Value2 = 7
$ws.Cells.Item(223,16).Value2 = 15
$ws.Cells.Item(223,17).Value2 = -2
$ws.Cells.Item(223,18).Value2 = 1.775
$ws.Cells.Item(223,19).Value2 = 2.025
$ws.Cells.Item(223,20).Value2 = 3.25
$ws.Cells.Item(223,21).Value2 = 1.825
$ws.Cells.Item(223,22).Value2 = 1.975
$ws.Cells.Item(223,23).Value2 = 0.1659999999999999
$ws.Cells.Item(223,24).Value2 = -1
$ws.Cells.Item(223,25).Value2 = -1
$ws.Cells.Item(223,26).Value2 = -1
$ws.Cells.Item(223,27).Value2 = 1.025
$ws.Cells.Item(223,28).Value2 = -1
$ws.Cells.Item(223,29).Value2 = 0.9750000000000001
$ws.Cells.Item(224,2).Value2 = 7623944
$ws.Cells.Item(224,5).Value2 = 45388.83333333334
$ws.Cells.Item(224,6).Value2 = 'Santos de Gupiles'
$ws.Cells.Item(224,7).Value2 = 'Municipal Liberia'
$ws.Cells.Item(224,11).Value2 = 2.9
$ws.Cells.Item(224,12).Value2 = 3.25
$ws.Cells.Item(224,13).Value2 = 2.375
$ws.Cells.Item(224,14).Value2 = 3
$ws.Cells.Item(224,15).Value2 = 3.3
$ws.Cells.Item(224,16).Value2 = 2.3
$ws.Cells.Item(224,17).Value2 = 0.25
$ws.Cells.Item(224,18).Value2 = 1.8
$ws.Cells.Item(224,19).Value2 = 2
$ws.Cells.Item(224,20).Value2 = 2.5
$ws.Cells.Item(224,21).Value2 = 1.85
$ws.Cells.Item(224,22).Value2 = 1.95
$ws.Cells.Item(225,1).Value2 = 223
$ws.Cells.Item(225,2).Value2 = 7623946
$ws.Cells.Item(225,3).Value2 = 'Costa Rica Primera Division'
$ws.Cells.Item(225,4).Value2 = 'Costa Rica Primera Division'
$ws.Cells.Item(225,5).Value2 = 45388.83333333334
$ws.Cells.Item(225,6).Value2 = 'Cartagines'
$ws.Cells.Item(225,7).Value2 = 'Alajuelense'
$ws.Cells.Item(225,11).Value2 = 3.4
$ws.Cells.Item(225,12).Value2 = 3.4
$ws.Cells.Item(225,13).Value2 = 1.95
$ws.Cells.Item(225,14).Value2 = 3.2
$ws.Cells.Item(225,15).Value2 = 3.4
$ws.Cells.Item(225,16).Value2 = 2.05
$ws.Cells.Item(225,17).Value2 = 0.25
$ws.Cells.Item(225,18).Value2 = 2
$ws.Cells.Item(225,19).Value2 = 1.8
$ws.Cells.Item(225,20).Value2 = 2.5
$ws.Cells.Item(225,21).Value2 = 1.825
$ws.Cells.Item(225,22).Value2 = 1.975
$ws.Cells.Item(225,23).Value2 = 0
$ws.Cells.Item(225,24).Value2 = 0
$ws.Cells.Item(225,25).Value2 = 0
$ws.Cells.Item(225,26).Value2 = 0
$ws.Cells.Item(225,27).Value2 = 0
$ws.Cells.Item(226,1).Value2 = 224
$ws.Cells.Item(226,2).Value2 = 7623996
$ws.Cells.Item(226,3).Value2 = 'Costa Rica Primera Division'
$ws.Cells.Item(226,4).Value2 = 'Costa Rica Primera Division'
$ws.Cells.Item(226,5).Value2 = 45388.92708333334
$ws.Cells.Item(226,6).Value2 = 'AD San Carlos'
$ws.Cells.Item(226,7).Value2 = 'Municipal Perez Zeledon'
$ws.Cells.Item(226,11).Value2 = 1.363
$ws.Cells.Item(226,12).Value2 = 4.5
$ws.Cells.Item(226,13).Value2 = 6.75
$ws.Cells.Item(226,14).Value2 = 1.363
$ws.Cells.Item(226,15).Value2 = 4.75
$ws.Cells.Item(226,16).Value2 = 6.5
$ws.Cells.Item(226,17).Value2 = -1.25
$ws.Cells.Item(226,18).Value2 = 1.8
$ws.Cells.Item(226,19).Value2 = 2
$ws.Cells.Item(226,20).Value2 = 3
$ws.Cells.Item(226,21).Value2 = 1.975
$ws.Cells.Item(226,22).Value2 = 1.825
$ws.Cells.Item(226,23).Value2 = 0
$ws.Cells.Item(226,24).Value2 = 0
$ws.Cells.Item(226,25).Value2 = 0
$ws.Cells.Item(226,26).Value2 = 0
$ws.Cells.Item(226,27).Value2 = 0
$ws.Cells.Item(227,1).Value2 = 225
$ws.Cells.Item(227,2).Value2 = 7623947
$ws.Cells.Item(227,3).Value2 = 'Costa Rica Primera Division'
$ws.Cells.Item(227,4).Value2 = 'Costa Rica Primera Division'
$ws.Cells.Item(227,5).Value2 = 45388.95833333334
$ws.Cells.Item(227,6).Value2 = 'Herediano'
$ws.Cells.Item(227,7).Value2 = 'AD Grecia'
$ws.Cells.Item(227,11).Value2 = 1.3
$ws.Cells.Item(227,12).Value2 = 4.75
$ws.Cells.Item(227,13).Value2 = 8
$ws.Cells.Item(227,14).Value2 = 1.333
$ws.Cells.Item(227,15).Value2 = 4.75
$ws.Cells.Item(227,16).Value2 = 7.5
$ws.Cells.Item(227,17).Value2 = -1.5
$ws.Cells.Item(227,18).Value2 = 2
$ws.Cells.Item(227,19).Value2 = 1.8
$ws.Cells.Item(227,20).Value2 = 2.75
$ws.Cells.Item(227,21).Value2 = 1.85
$ws.Cells.Item(227,22).Value2 = 1.95
$ws.Cells.Item(227,23).Value2 = 0
$ws.Cells.Item(227,24).Value2 = 0
$ws.Cells.Item(227,25).Value2 = 0
$ws.Cells.Item(227,26).Value2 = 0
$ws.Cells.Item(227,27).Value2 = 0
$ws.Cells.Item(228,1).Value2 = 226
$ws.Cells.Item(228,2).Value2 = 8048492
$ws.Cells.Item(228,3).Value2 = 'Costa Rica Primera Division'
$ws.Cells.Item(228,4).Value2 = 'Costa Rica Primera Division'
$ws.Cells.Item(228,5).Value2 = 45389.75
$ws.Cells.Item(228,6).Value2 = 'AD Guanacasteca'
$ws.Cells.Item(228,7).Value2 = 'Puntarenas'
$ws.Cells.Item(228,11).Value2 = 2.2
$ws.Cells.Item(228,12).Value2 = 3.2
$ws.Cells.Item(228,13).Value2 = 3.3
$ws.Cells.Item(228,14).Value2 = 2.15
$ws.Cells.Item(228,15).Value2 = 3.2
$ws.Cells.Item(228,16).Value2 = 3.5
$ws.Cells.Item(228,17).Value2 = -0.25
$ws.Cells.Item(228,18).Value2 = 1.825
$ws.Cells.Item(228,19).Value2 = 1.975
$ws.Cells.Item(228,20).Value2 = 2.25
$ws.Cells.Item(228,21).Value2 = 1.975
$ws.Cells.Item(228,22).Value2 = 1.825
$ws.Cells.Item(228,23).Value2 = 0
$ws.Cells.Item(228,24).Value2 = 0
$ws.Cells.Item(228,25).Value2 = 0
$ws.Cells.Item(228,26).Value2 = 0
$ws.Cells.Item(228,27).Value2 = 0
$ws.Cells.Item(229,1).Value2 = 227
$ws.Cells.Item(229,2).Value2 = 7623997
$ws.Cells.Item(229,3).Value2 = 'Costa Rica Primera Division'
$ws.Cells.Item(229,4).Value2 = 'Costa Rica Primera Division'
$ws.Cells.Item(229,5).Value2 = 45389.83333333334
$ws.Cells.Item(229,6).Value2 = 'Deportivo Saprissa'
$ws.Cells.Item(229,7).Value2 = 'Sporting San Jose'
$ws.Cells.Item(229,11).Value2 = 1.285
$ws.Cells.Item(229,12).Value2 = 5
$ws.Cells.Item(229,13).Value2 = 8
$ws.Cells.Item(229,14).Value2 = 1.3
$ws.Cells.Item(229,15).Value2 = 5
$ws.Cells.Item(229,16).Value2 = 7.5
$ws.Cells.Item(229,17).Value2 = -1.5
$ws.Cells.Item(229,18).Value2 = 1.975
$ws.Cells.Item(229,19).Value2 = 1.825
$ws.Cells.Item(229,20).Value2 = 2.75
$ws.Cells.Item(229,21).Value2 = 1.875
$ws.Cells.Item(229,22).Value2 = 1.925
$ws.Cells.Item(229,23).Value2 = 0
$ws.Cells.Item(229,24).Value2 = 0
$ws.Cells.Item(229,25).Value2 = 0
$ws.Cells.Item(229,26).Value2 = 0
$ws.Cells.Item(229,27).Value2 = 0
